$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.402.82'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').Value = '3.593.80'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.50'
$ws.Range('E5').Value = '  +2.40%  '
$ws.Range('E6').Value = '  +18.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '653.46'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.425'
$ws.Range('E8').Value = '  +5.68%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.07'
$ws.Range('E9').Value = '  +3.71%  '
$ws.Range('B10').Value = 'USDC'
$ws.Range('C10').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.999'
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').Value = '3.590.60'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.03'
$ws.Range('E12').Value = '  +4.99%  '
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '4.261.09'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '97.253.49'
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000260'
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('D18').Value = '3.582.33'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.75'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.63'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.28'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.536'
$ws.Range('E22').Value = '  +8.54%  '
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '517.73'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('E25').Value = '  +3.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.98'
$ws.Range('E26').Value = '  -1.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '103.20'
$ws.Range('E27').Value = '  +7.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.33'
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.182'
$ws.Range('E29').Value = '  +25.17%  '
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.05'
$ws.Range('E31').Value = '  +4.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('E33').Value = '  +6.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '31.84'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E37').Value = '  +3.75%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '618.17'
$ws.Range('E38').Value = '  +3.47%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.80'
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.451'
$ws.Range('E44').Value = '  +37.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.12'
$ws.Range('E45').Value = '  +5.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0453'
$ws.Range('E46').Value = '  +8.50%  '
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.67'
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.75'
$ws.Range('E49').Value = '  +6.13%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.29'
$ws.Range('E50').Value = '  +6.63%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '32.80'
$ws.Range('E51').Value = '  -4.44%  '
